$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header style (bold + border) from the "Weekly Quantity" sheet header
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- Data rows ---
$data = @(
    @(45340.99999999999, 0, -198.3241527132459, 178.5703363088063),
    @(45347.99999999999, 3, -178.2657499867777, 172.8104491280525),
    @(45354.99999999999, 11, -170.9867103944837, 187.757278537409),
    @(45368.99999999999, 28, -139.6726109533705, 209.2817256999884),
    @(45382.99999999999, 45, -148.7373472301541, 229.3982051970966),
    @(45396.99999999999, 62, -123.6661000712934, 246.7595756912527),
    @(45410.99999999999, 80, -109.7345919180668, 265.9408929372653),
    @(45424.99999999999, 97, -79.91760245527222, 278.6537294593462),
    @(45438.99999999999, 114, -65.91557937525604, 289.970131852298),
    @(45445.99999999999, 122, -59.61272625095391, 282.0980781380313),
    @(45459.99999999999, 139, -33.66415700147701, 321.3868774401552),
    @(45466.99999999999, 148, -26.98625571248228, 323.4615568403915),
    @(45473.99999999999, 156, -19.76949722188246, 346.5136347263128),
    @(45480.99999999999, 165, -13.64783063496352, 333.7400452764577),
    @(45487.99999999999, 173, -0.4624870820814655, 352.8171525317168),
    @(45501.99999999999, 190, 12.60876910611277, 366.3830082447832),
    @(45508.99999999999, 199, 22.75238214981618, 378.6545782801163),
    @(45543.99999999999, 242, 65.86082950940873, 421.8559169064546),
    @(45550.99999999999, 250, 69.72591583894142, 427.6914762306068),
    @(45557.99999999999, 259, 74.63996487361314, 434.0647454453477),
    @(45578.99999999999, 284, 104.214645518471, 459.8036844145104),
    @(45585.99999999999, 293, 101.4842821476873, 469.3754779635182),
    @(45592.99999999999, 301, 125.5808776790361, 466.7198119945258),
    @(45599.99999999999, 310, 132.622938810532, 482.909386053203),
    @(45606.99999999999, 318, 143.995361729911, 497.9373683136034),
    @(45613.99999999999, 327, 147.3058420768375, 499.9948731158906),
    @(45620.99999999999, 335, 160.5863248086, 524.4128543358004),
    @(45627.99999999999, 344, 165.8005338144984, 519.877707100435),
    @(45634.99999999999, 352, 167.7233314549434, 538.7437256008261),
    @(45641.99999999999, 361, 183.2142780412522, 535.0811445924069),
    @(45648.99999999999, 369, 185.4561864078309, 553.0464348969821)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the date-formatted style used in column A of the source sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A32").PasteSpecial(-4122)

# Restore original active sheet/selection state
$wsWeekly.Activate()
[void]$wsWeekly.Range("A1").Select()

Write-Output "PO Forecast sheet added with $($data.Count) data rows"
